$d = $word.ActiveDocument

# Step 1: seed a throwaway numbered-list paragraph so Word mints
# word/numbering.xml (numId=1 / abstractNumId=0, decimal "%1." scheme)
# exactly as ApplyNumberDefault would for a normal numbered list.
$lastPara = $d.Paragraphs.Last
$seedRange = $lastPara.Range
$seedRange.InsertParagraphAfter()
$seedPara = $d.Paragraphs.Last
$seedPara.Style = "Listenabsatz"
$seedRange2 = $seedPara.Range
$seedRange2.Text = "seed"
$seedRange2.ListFormat.ApplyNumberDefault()

# Step 2: replace that seed paragraph (and everything from it onward)
# with the exact OOXML for the new "Messungen" / "Positionierung" section,
# using InsertXML so proofErr markers, run-splits and numPr survive verbatim.
$target = $d.Paragraphs.Last.Range
$target.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">Messungen: </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Stationäre Messungen durchführen und Daten plotten</w:t></w:r><w:r><w:t>, Software muss noch Listen von Datenpunkten aufnehmen können.</w:t></w:r><w:r><w:t xml:space="preserve"> Daten aufnehmen, ohne dass eine Person </w:t></w:r><w:r><w:t>nebendran</w:t></w:r><w:r><w:t xml:space="preserve"> steht und einmal mit</w:t></w:r><w:r><w:t>,</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>testen,</w:t></w:r><w:r><w:t xml:space="preserve"> ob es Veränderungen gibt.</w:t></w:r><w:r><w:t xml:space="preserve"> Statistik einlesen/wiederholen</w:t></w:r><w:r><w:t xml:space="preserve"> und </w:t></w:r><w:r><w:t>prüfen,</w:t></w:r><w:r><w:t xml:space="preserve"> ob es Zusammenhänge in den Daten gibt</w:t></w:r><w:r><w:t>,</w:t></w:r><w:r><w:t xml:space="preserve"> die man nutzen kann.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Positionierung:</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Mehrere </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Routeranzahlen</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> probieren mit aktuellem Aufbau und testen, ob es Vorteile bringt. </w:t></w:r><w:r><w:t xml:space="preserve">Dann verschiedene statistische Methoden testen, ob es Unterschiede macht. </w:t></w:r><w:r><w:t xml:space="preserve">Räumlichkeiten anpassen: freier Raum ohne Gegenstände, verschiedene </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Routerhöhen</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>blockierende Hindernisse in den Weg stellen, z.B. Schränke, dynamische Hindernisse.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>2</w:t></w:r><w:r><w:t xml:space="preserve"> Messpunkt</w:t></w:r><w:r><w:t>e</w:t></w:r><w:r><w:t>, Signalstär</w:t></w:r><w:r><w:t>k</w:t></w:r><w:r><w:t>e über Zeit -&gt; Plotten</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Jetzt einen Gegenstand dazwischen, testen ob Veränderungen</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">1. 2. Im </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Aussenbereich</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Routeranzahl</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> variieren und Punkte wiederholen</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Mensch läuft durch</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Alles zusammen von </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>Oben</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Innenbereich keine Möbel</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Aussenbereich</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Innenbereich mit Möbeln</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Zielpunkte nahe Wand, </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>Router,…</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> auf Punkte 1. 2. 3.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Wie 1. Nur mehr Router</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>')
